$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.951.49"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").Value = "3.512.42"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.91%  "
$ws.Range("D7").Value = "3.513.62"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.378"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.49%  "
$ws.Range("D13").Value = "4.115.54"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.539.24"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000178"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("D18").Value = "64.003.45"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("E19").Value = "  -3.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "384.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.574"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("D24").Value = "3.654.83"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.69%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000115"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("E32").Value = "  -1.71%  "
$ws.Range("D33").Value = "3.522.66"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.71%  "
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "159.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0789"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.813"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.62%  "
$ws.Range("E46").Value = "  -5.23%  "
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.82%  "
$ws.Range("D49").Value = "2.439.61"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.910"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.09%  "
